# Financials update: insert a new "latest period" column before column D,
# shifting the existing 8 periods of data from D:K to E:L, then populate the
# new column D with the newest period's figures (plus a handful of small
# corrections to previously-reported figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D. Excel shifts D:K -> E:L automatically,
# including formulas/values, carrying their original styles with them.
$ws.Columns("D:D").Insert()

# The brand-new column D cells default to the column's base style; copy the
# real per-row number formats/styles over from column E (which just received
# the original column D's formatting) so D matches the rest of each row.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Populate column D with the newest reporting period ----

# Period Ending headers (D7/D38/D80)
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# Income Statement (rows 8-35)
$ws.Range("D8").Value = 6426000
$ws.Range("D9").Value = 3829000
$ws.Range("D10").Value = 2597000
$ws.Range("D12").Value = 2000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 20000
$ws.Range("D15").Value = 424000
$ws.Range("D17").Value = 4459000
$ws.Range("D18").Value = 1967000
$ws.Range("D20").Value = 40000
$ws.Range("D21").Value = 2431000
$ws.Range("D22").Value = 733000
$ws.Range("D23").Value = 1274000
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1274000
$ws.Range("D27").Value = 1214000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -40000
$ws.Range("D33").Value = 1214000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1214000

# Balance Sheet - Assets (rows 41-54)
$ws.Range("D41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 690000
$ws.Range("D44").Value = 99000
$ws.Range("D45").Value = 1567000
$ws.Range("D46").Value = 2356000
$ws.Range("D47").Value = 157000
$ws.Range("D48").Value = 15390000
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 71000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 17974000

# Balance Sheet - Liabilities (rows 57-77)
$ws.Range("D57").Value = 15000
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 1053000
$ws.Range("D60").Value = 1068000
$ws.Range("D61").Value = 16066000
$ws.Range("D62").Value = 40000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 17174000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = -16000
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 816000
$ws.Range("D77").Value = 0

# Cash Flow (rows 81-102)
$ws.Range("D81").Value = 1214000
$ws.Range("D83").Value = 424000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1874000
$ws.Range("D91").Value = -804000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -804000
$ws.Range("D96").Value = -1113000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1118000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -48000

# ---- A few small corrections to previously reported periods (now in ----
# ---- their shifted columns) that are not pure carry-overs.          ----

# Row 43 (Goodwill): old I43 (now J43) corrected from 7000 -> 7200
$ws.Range("J43").Value = 7200

# Row 44 (Intangible Assets): old I44 (now J44) corrected from 7000 -> 9700
$ws.Range("J44").Value = 9700

# Row 47 (Deferred Long Term Asset Charges): the six older zero periods are
# recast as "NA", and the previously-blank last column now reports 0
$ws.Range("E47:J47").Value = "NA"
$ws.Range("K47").Value = 0

# Row 59 (Short/Current Long Term Debt): old I59 (now J59) corrected from
# 88600 -> 88500
$ws.Range("J59").Value = 88500
